$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number that was updated
# from 45204 (2023-10-05) to 45207 (2023-10-08) for every data row (2-32).
$range = $ws.Range("C2:C32")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 45204) {
        $cell.Value2 = 45207
    }
}
